$wb = $excel.ActiveWorkbook

# Rename sheets to reflect the deeper chapter numbering (1.1.A/B/C -> 1.1.1.A/B/C)
$wb.Worksheets.Item("1.1.A").Name = "1.1.1.A"
$wb.Worksheets.Item("1.1.B").Name = "1.1.1.B"
$wb.Worksheets.Item("1.1.C").Name = "1.1.1.C"

# Activate the completed/tested sheet (1.1.1.C) and put the selection at B3
$ws = $wb.Worksheets.Item("1.1.1.C")
$ws.Activate()
$ws.Range("B3").Select()

Write-Output "done"
